$wb = $excel.ActiveWorkbook

# ALC row 106: Enchanted Palladium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 148151460
$ws.Range("I106").Value = 66671132
$ws.Range("J106").Value = 250001870
$ws.Range("K106").Value = 66671132
$ws.Range("L106").Value = 250001870
$ws.Range("M106").Value = -66670501
$ws.Range("N106").Value = -250003132

# ALC row 137: Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1112.0857
$ws.Range("I137").Value = 919.4355
$ws.Range("J137").Value = 2605.125
$ws.Range("K137").Value = 2758.3065
$ws.Range("L137").Value = 7815.375
$ws.Range("M137").Value = -208.3065000000001
$ws.Range("N137").Value = -12915.375

# ALC row 138: Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2022.1522
$ws.Range("I138").Value = 844.3158
$ws.Range("J138").Value = 3940.3428
$ws.Range("K138").Value = 2532.9474
$ws.Range("L138").Value = 11821.0284
$ws.Range("M138").Value = 2607.0526
$ws.Range("N138").Value = -22101.0284

# ARM row 32: Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1726.15
$ws.Range("I32").Value = 1659.8632
$ws.Range("J32").Value = 2985.6
$ws.Range("K32").Value = 1659.8632
$ws.Range("L32").Value = 2985.6
$ws.Range("M32").Value = -1372.8632
$ws.Range("N32").Value = -3559.6

# ARM row 55: Mythril Elmo
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 18990
$ws.Range("J55").Value = 18990
$ws.Range("L55").Value = 18990
$ws.Range("N55").Value = -19620

# ARM row 63: Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 125003090
$ws.Range("I63").Value = 166670420
$ws.Range("J63").Value = 1115
$ws.Range("K63").Value = 166670420
$ws.Range("L63").Value = 1115
$ws.Range("M63").Value = -166669734
$ws.Range("N63").Value = -2487

# ARM row 66: Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 125003090
$ws.Range("I66").Value = 166670420
$ws.Range("J66").Value = 1115
$ws.Range("K66").Value = 833352100
$ws.Range("L66").Value = 5575
$ws.Range("M66").Value = -833348668
$ws.Range("N66").Value = -12439

# ARM row 74: Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1374.3448
$ws.Range("I74").Value = 1228.5454
$ws.Range("J74").Value = 1832.5714
$ws.Range("K74").Value = 1228.5454
$ws.Range("L74").Value = 1832.5714
$ws.Range("M74").Value = -354.5454
$ws.Range("N74").Value = -3580.5714

# ARM row 77: Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1374.3448
$ws.Range("I77").Value = 1228.5454
$ws.Range("J77").Value = 1832.5714
$ws.Range("K77").Value = 6142.727
$ws.Range("L77").Value = 9162.857
$ws.Range("M77").Value = -1774.727
$ws.Range("N77").Value = -17898.857

# ARM row 80: Titanium Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 23329
$ws.Range("I80").Value = 13938.5
$ws.Range("J80").Value = 42110
$ws.Range("K80").Value = 13938.5
$ws.Range("L80").Value = 42110
$ws.Range("M80").Value = -12940.5
$ws.Range("N80").Value = -44106

# ARM row 83: Titanium Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 23329
$ws.Range("I83").Value = 13938.5
$ws.Range("J83").Value = 42110
$ws.Range("K83").Value = 41815.5
$ws.Range("L83").Value = 126330
$ws.Range("M83").Value = -36823.5
$ws.Range("N83").Value = -136314

# BSM row 94: High Steel Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1146.0555
$ws.Range("I94").Value = 433
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 433
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = 18
$ws.Range("N94").Value = -3902

# CRP row 31: Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6000.6284
$ws.Range("I31").Value = 1471.5319
$ws.Range("J31").Value = 15255.739
$ws.Range("K31").Value = 1471.5319
$ws.Range("L31").Value = 15255.739
$ws.Range("M31").Value = -1176.5319
$ws.Range("N31").Value = -15845.739

# CRP row 34: Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6000.6284
$ws.Range("I34").Value = 1471.5319
$ws.Range("J34").Value = 15255.739
$ws.Range("K34").Value = 1471.5319
$ws.Range("L34").Value = 15255.739
$ws.Range("M34").Value = -1269.5319
$ws.Range("N34").Value = -15659.739

# CUL row 75: Emerald Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 23811088
$ws.Range("J75").Value = 23811088
$ws.Range("L75").Value = 71433264
$ws.Range("N75").Value = -71435260

# CUL row 78: Emerald Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 23811088
$ws.Range("J78").Value = 23811088
$ws.Range("L78").Value = 214299792
$ws.Range("N78").Value = -214309776

# CUL row 87: Clam Chowder
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 5050
$ws.Range("I87").Value = 5133.3335
$ws.Range("K87").Value = 15400.0005
$ws.Range("M87").Value = -14152.0005

# CUL row 90: Clam Chowder
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 5050
$ws.Range("I90").Value = 5133.3335
$ws.Range("K90").Value = 46200.0015
$ws.Range("M90").Value = -39960.0015

# CUL row 108: Grilled Rail
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 33644.5
$ws.Range("I108").Value = 33644.5
$ws.Range("K108").Value = 100933.5
$ws.Range("M108").Value = -98053.5

# CUL row 113: Night Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 179087.5
$ws.Range("I113").Value = 519.2059
$ws.Range("J113").Value = 455056.7
$ws.Range("K113").Value = 1557.6177
$ws.Range("L113").Value = 1365170.1
$ws.Range("M113").Value = 612.3822999999998
$ws.Range("N113").Value = -1369510.1

# CUL row 121: Coffee Biscuit
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 963.1613
$ws.Range("I121").Value = 496
$ws.Range("J121").Value = 1053
$ws.Range("K121").Value = 1488
$ws.Range("L121").Value = 3159
$ws.Range("M121").Value = -178
$ws.Range("N121").Value = -5779

# CUL row 122: Northern Sea Salt
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2811.756
$ws.Range("I122").Value = 450.75
$ws.Range("J122").Value = 4322.8
$ws.Range("K122").Value = 4056.75
$ws.Range("L122").Value = 38905.2
$ws.Range("M122").Value = -1606.75
$ws.Range("N122").Value = -43805.2

# GSM row 126: Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5118.4414
$ws.Range("I126").Value = 11571.2
$ws.Range("J126").Value = 2429.7917
$ws.Range("K126").Value = 34713.60000000001
$ws.Range("L126").Value = 7289.375100000001
$ws.Range("M126").Value = -32243.60000000001
$ws.Range("N126").Value = -12229.3751

# GSM row 132: Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1675.738
$ws.Range("I132").Value = 1288.3611
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 3865.0833
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -1335.0833
$ws.Range("N132").Value = -17060

# LTW row 133: Loboskin Amulet of Fending
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# LTW row 136: Br'aax Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4554.607
$ws.Range("I136").Value = 3946.8865
$ws.Range("J136").Value = 6782.9165
$ws.Range("K136").Value = 11840.6595
$ws.Range("L136").Value = 20348.7495
$ws.Range("M136").Value = -9290.6595
$ws.Range("N136").Value = -25448.7495

# WVR row 132: Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19394.834
$ws.Range("I132").Value = 25662.824
$ws.Range("J132").Value = 1486.2858
$ws.Range("K132").Value = 76988.47200000001
$ws.Range("L132").Value = 4458.857400000001
$ws.Range("M132").Value = -74458.47200000001
$ws.Range("N132").Value = -9518.857400000001
